$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 11;  D = "[1, 0, 0, 0, 1, 0, 0]"; E = "['Normal', 'RegulationViolation']" },
    @{ Row = 25;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 26;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 29;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 38;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 58;  D = "[0, 0, 0, 1, 0, 0, 0]"; E = "['ParamViolation']" },
    @{ Row = 67;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 75;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 82;  D = "[1, 1, 0, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment']" },
    @{ Row = 88;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 92;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 107; D = "[1, 0, 0, 0, 0, 1, 0]"; E = "['Normal', 'CommunicationIssue']" },
    @{ Row = 113; D = "[1, 0, 1, 0, 0, 0, 1]"; E = "['Normal', 'HardwareFault', 'SoftwareFault']" },
    @{ Row = 116; D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" }
)

foreach ($u in $updates) {
    $ws.Range("D" + $u.Row).Value = $u.D
    $ws.Range("E" + $u.Row).Value = $u.E
}

$wb.Save()
